$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new rows 73-75 (id/en first, then th) -------------------------------
$ws.Range("A73").Value = "youCanCheckOrderTxt2"
$ws.Range("B73").Value = "by input your order no."
$ws.Range("A74").Value = "youCanCheckOrderTxt3"
$ws.Range("B74").Value = "and your email"

$ws.Range("C73").Value = "โดยใส่เลขคำสั่งซื้อคือ"
$ws.Range("C74").Value = "และอีเมล์ของท่าน"

$ws.Range("A75").Value = "errFailToUploadQr"
$ws.Range("B75").Value = "Fail to upload QR"
$ws.Range("C75").Value = "ไม่สามารถสร้างคิวอาร์ได้"

# --- Fix existing error message rows 67-68 ------------------------------------
$ws.Range("C67").Value = "ไม่พบหมายเลขคำสั่งซื้อ"
$ws.Range("B67").Value = "Invalid Order Id"

$ws.Range("B68").Value = "No file uploaded"
$ws.Range("C68").Value = "ไม่พบรูปหลักฐานการโอนเงิน"

# --- Add new row 76 ------------------------------------------------------------
$ws.Range("A76").Value = "errFailToUploadEvidence"
$ws.Range("B76").Value = "Fail to upload evidence"
$ws.Range("C76").Value = "ไม่สามารถส่งหลักฐานได้"

# --- Fill the remaining locale columns (D:I) with the English fallback value --
$ws.Range("D73:I73").Value = "by input your order no."
$ws.Range("D74:I74").Value = "and your email"
$ws.Range("D75:I75").Value = "Fail to upload QR"
$ws.Range("D76:I76").Value = "Fail to upload evidence"

$ws.Range("D67:I67").Value = "Invalid Order Id"

# --- Update the view state to match the author's final selection --------------
$ws.Range("D76:I76").Select()
$excel.ActiveWindow.ScrollRow = 56
